{"js": "// The edit removes the trailing \"Comments graph:\" review/feedback section\n// (three blank paragraphs, a page break, the \"Comments graph:\" heading, an\n// empty paragraph, and the bulleted list of comments) that used to follow\n// the \"Annual budget and multi-year financial plan\" paragraph at the very\n// end of the document body - right before the section properties.\n//\n// We locate the anchor paragraph by its distinctive trailing sentence so\n// the script is resilient to absolute paragraph-index changes, then remove\n// every paragraph that follows it through the end of the document body.\n\nconst body = context.document.body;\n\n// 1) Find the paragraph that ends the section we want to keep.\nconst anchorText = \"This integration supports fiscal discipline, transparency, and sustainable resource management.\";\nconst searchResults = body.search(anchorText, { matchCase: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Anchor paragraph not found - document may already be edited.\");\n}\n\nconst anchorRange = searchResults.items[0];\nconst anchorParagraph = anchorRange.paragraphs.getFirst();\n\n// 2) Collect every paragraph after the anchor paragraph (these are the\n// blank spacer paragraphs, the page break, \"Comments graph:\" and the\n// bulleted review comments) and delete them one at a time, starting from\n// the paragraph right after the anchor and always re-fetching \"the next\n// paragraph\" so each delete call only ever targets a paragraph that still\n// has a sibling after it (the very last paragraph of a section cannot be\n// removed while it is empty, but removing paragraphs in this order never\n// leaves an empty paragraph stranded as the final one).\nlet next = anchorParagraph.getNextOrNullObject();\nnext.load(\"isNullObject\");\nawait context.sync();\n\nwhile (!next.isNullObject) {\n  next.delete();\n  await context.sync();\n\n  next = anchorParagraph.getNextOrNullObject();\n  next.load(\"isNullObject\");\n  await context.sync();\n}\n", "ps1": "# The edit removes the trailing \"Comments graph:\" review/feedback section\n# (three blank paragraphs, a page break, the \"Comments graph:\" heading, an\n# empty paragraph, and the bulleted list of comments) that used to follow\n# the \"Annual budget and multi-year financial plan\" paragraph at the very\n# end of the document body - right before the section properties.\n#\n# The anchor paragraph is located by its distinctive trailing sentence so\n# the script is resilient to absolute paragraph-index changes, then every\n# paragraph that follows it through the end of the document body is\n# removed.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"sustainable resource management.\"\n\n$paragraphs = $d.Paragraphs\n$anchorIndex = -1\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    if ($paragraphs.Item($i).Range.Text -like (\"*\" + $anchorText + \"*\")) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph not found - document may already be edited.\"\n}\n\n# Repeatedly delete the paragraph right after the anchor paragraph. Doing\n# it this way (instead of deleting from the end backwards) means a delete\n# call never targets the very last, now-empty paragraph of the section\n# while it is still the sole remaining paragraph of that range, which this\n# host does not allow.\nwhile ($d.Paragraphs.Count -gt $anchorIndex) {\n    $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n}\n"}
